# Petty cash book update - 7-Jan-2021, midday update.
# Adds new transaction rows (4-9, 11-16, 18-27, 29-33) for 3 additional days
# (2021-01-01 .. 2021-01-03 in Excel serial form: 44201, 44202, 44203) plus a
# new Debit entry on the existing first day (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# --- 2021-01-07 (existing day, row 3) : add another Wages Expense debit ---
$ws.Range("D3").Formula = "=45000+240000"

# --- Day 44200 continued : rows 4-9 ---
$ws.Range("B4").Value = "TRANSFER BCA"
$ws.Range("D4").Formula = "=5170000+1040000+450000+7374000+1338000+586000+3467000+3219000+5170000+40000000+920000"

$ws.Range("B5").Value = "A/R"
$ws.Range("C5").Formula = "=7374000+2211000+1500000+100000000+94174000"

$ws.Range("B6").Value = "BELI sanitizer"
$ws.Range("D6").Formula = "=50000"

$ws.Range("B7").Value = "SALES - cash/retail"
$ws.Range("C7").Formula = "=166645025-43610025-94174000"

$ws.Range("B8").Value = "SELISIH - lebih"
$ws.Range("C8").Value = 70000

$ws.Range("B9").Value = "SETOR KE BANK"
$ws.Range("D9").Value = 165000000

# --- Day 44201 : rows 10-16 ---
$ws.Range("A10").Value = 44201
$ws.Range("B10").Value = "Wages Expense"
$ws.Range("D10").Formula = "=45000+150000+195000"

$ws.Range("B11").Value = "TRANSFER BCA"
$ws.Range("D11").Formula = "=5540000+1658000+1012500+440000"

$ws.Range("B12").Value = "FREIGHT-IN"
$ws.Range("D12").Value = 1600000

$ws.Range("B13").Value = "SALES - cash/retail"
$ws.Range("C13").Formula = "=18006025+9890475-20113000"

$ws.Range("B14").Value = "A/R"
$ws.Range("C14").Formula = "=20113000"

$ws.Range("B15").Value = "SELISIH - lebih"
$ws.Range("C15").Value = 20000

$ws.Range("B16").Value = "SETOR KE BANK"
$ws.Range("D16").Value = 17000000

# --- Day 44202 : rows 17-27 ---
$ws.Range("A17").Value = 44202
$ws.Range("B17").Value = "Wages Expense"
$ws.Range("D17").Formula = "=45000+255000"

$ws.Range("B18").Value = "TRANSFER BCA"
$ws.Range("D18").Formula = "=680000+22000000+6382000+1227000"

$ws.Range("B19").Value = "A/R"
$ws.Range("C19").Formula = "=2000000+5000000+100000000+82801000"

$ws.Range("B20").Value = "BELI kresek"
$ws.Range("D20").Formula = "=95000"

$ws.Range("B21").Value = "CHEQUE RECEIVED"
$ws.Range("D21").Formula = "=2169000"

$ws.Range("B22").Value = "IURAN DAERAH"
$ws.Range("D22").Formula = "=25000"

$ws.Range("B23").Value = "A/P"
$ws.Range("D23").Formula = "=1266000"

$ws.Range("B24").Value = "BENSIN - RUSH"
$ws.Range("D24").Value = 250000

$ws.Range("B25").Value = "SALES - cash/retail"
$ws.Range("C25").Formula = "=166021025-73732025-82801000"

$ws.Range("B26").Value = "SELISIH - lebih"
$ws.Range("C26").Value = 68000

$ws.Range("B27").Value = "SETOR KE BANK"
$ws.Range("D27").Value = 165000000

# --- Day 44203 : rows 28-33 ---
$ws.Range("A28").Value = 44203
$ws.Range("B28").Value = "Wages Expense"
$ws.Range("D28").Formula = "=45000"

$ws.Range("B29").Value = "A/R"
$ws.Range("C29").Formula = "=19800000+851000+9149000"

$ws.Range("B30").Value = "TRANSFER DANAMON"
$ws.Range("D30").Formula = "=19800000"

$ws.Range("B31").Value = "PARKIR - bulanan"
$ws.Range("D31").Value = 10000

$ws.Range("B32").Value = "TRANSFER BCA"
$ws.Range("D32").Formula = "=215000+308000+1240000"

$ws.Range("B33").Value = "A/P"
$ws.Range("D33").Value = 6400000

# --- Scroll the frozen pane down so the newest entries are in view, and
#     leave the active selection on C24 (matches the author's cursor spot) ---
$ws.Range("C24").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
